$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price & volume data (GitHub Actions refresh)

$ws.Range("D2").Value = "26.230.84"
$ws.Range("E2").Value = "  -2.03%  "

$ws.Range("D3").Value = "1.673.81"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.78"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5141"
$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2642"
$ws.Range("E8").Value = "  +1.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06417"
$ws.Range("E9").Value = "  +4.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07416"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").Value = "1.674.62"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5827"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008655"
$ws.Range("E15").Value = "  +6.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.55"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").Value = "26.278.68"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.955"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.65"
$ws.Range("E21").Value = "  +3.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.217"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.84"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.650"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1185"
$ws.Range("E26").Value = "  +3.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.67"
$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06424"
$ws.Range("E28").Value = "  +13.73%  "

$ws.Range("E29").Value = "  -1.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.326"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.525"
$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.521"
$ws.Range("E32").Value = "  +2.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.637"
$ws.Range("E33").Value = "  -2.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.021"
$ws.Range("E34").Value = "  +1.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6067"
$ws.Range("E35").Value = "  +2.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.375"
$ws.Range("E36").Value = "  -1.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.666"
$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.202"
$ws.Range("E38").Value = "  +4.16%  "

$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").Value = "1.088.30"
$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8659"
$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.32"
$ws.Range("E43").Value = "  +3.11%  "

$ws.Range("D44").Value = "1.823.73"
$ws.Range("E44").Value = "  -1.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("E45").Value = "  +6.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.34"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.124"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05217"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4303"
$ws.Range("E50").Value = "  -0.81%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.897"
$ws.Range("E51").Value = "  +3.41%  "
